# feat: add 2022-Q1 data
#
# The workbook originally has two sheets: "2021-Q4" (holdings detail) and
# "总计" (summary totals, one row per quarter). This change inserts a new
# "2022-Q1" holdings sheet (same shape/format as "2021-Q4") between the two
# existing sheets, and adds a matching "2022-Q1" summary row at the top of
# the "总计" sheet (pushing the existing "2021-Q4" row down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" so that it
#    starts out with identical header text/styles/column layout, then
#    place it right after "2021-Q4" (i.e. before "总计").
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.ActiveSheet
$q1.Name = "2022-Q1"

# Helper: write a value as TEXT (shared string) without leaving a residual
# "@ text" number-format style behind - matches how the source file stores
# numeric-looking strings (fund codes, percentages, etc.) as plain text
# cells with no explicit style.
function Set-TextCell($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Fund holdings data for 2022-Q1 (index, code, name, scale, totalPosition,
# positionPct, marketValue, positionRank).
$rows = @(
    @(0,  "320003", "诺安先锋混合",                 "45.79", "69.96", "3.93", "1.7995", 5),
    @(1,  "001743", "诺安优选回报灵活配置混合",       "6.13",  "71.32", "8.98", "0.5505", 1),
    @(2,  "688888", "浙商聚潮产业成长混合",           "8.25",  "93.40", "4.63", "0.3820", 9),
    @(3,  "320018", "诺安新动力混合",                 "3.36",  "79.24", "2.90", "0.0974", 10),
    @(4,  "009568", "浙商智多宝稳健一年持有期混合A",   "3.13",  "22.64", "0.97", "0.0304", 8),
    @(5,  "007439", "东海科技动力混合A",              "0.36",  "90.01", "6.40", "0.0230", 1),
    @(6,  "009569", "浙商智多宝稳健一年持有期混合C",   "1.59",  "22.64", "0.97", "0.0154", 8),
    @(7,  "007463", "东海科技动力混合C",              "0.20",  "90.01", "6.40", "0.0128", 1),
    @(8,  "005901", "诺安汇利灵活配置混合A",           "0.08",  "86.88", "6.02", "0.0048", 6),
    @(9,  "001744", "诺安进取回报灵活配置混合",        "0.04",  "62.10", "5.89", "0.0024", 1),
    @(10, "006538", "东海核心价值精选混合",            "0.03",  "89.24", "6.36", "0.0019", 3),
    @(11, "005902", "诺安汇利灵活配置混合C",           "0.02",  "86.88", "6.02", "0.0012", 6)
)

$lastExistingRow = $q4.UsedRange.Rows.Count   # 8 (header + 7 existing data rows)
$neededRows = $rows.Count                      # 12 data rows needed
$lastNeededRow = 1 + $neededRows               # row 13

# Extend the sheet with extra (empty but correctly-styled) rows by copying
# the format of the last existing data row down as far as we need.
if ($lastNeededRow -gt $lastExistingRow) {
    $srcRow = "A" + $lastExistingRow + ":H" + $lastExistingRow
    $dstRange = "A" + ($lastExistingRow + 1) + ":H" + $lastNeededRow
    $q1.Range($srcRow).Copy()
    $q1.Range($dstRange).PasteSpecial(-4122)   # xlPasteFormats
}

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q1.Range("A$r").Value = $row[0]
    Set-TextCell $q1.Range("B$r") $row[1]
    Set-TextCell $q1.Range("C$r") $row[2]
    Set-TextCell $q1.Range("D$r") $row[3]
    Set-TextCell $q1.Range("E$r") $row[4]
    Set-TextCell $q1.Range("F$r") $row[5]
    Set-TextCell $q1.Range("G$r") $row[6]
    $q1.Range("H$r").Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing "2021-Q4" row
#    down to row 3 and insert the new "2022-Q1" totals in row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give row 3's index cell (A3) the same styling as A2 before overwriting.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
Set-TextCell $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.78

$total.Range("A2").Value = 0
Set-TextCell $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 2.92
